# Update Leve profit/price figures pulled from the market-board scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 6288.7744
$ws.Range("I76").Value = 7150.0835
$ws.Range("J76").Value = 3335.7144
$ws.Range("K76").Value = 7150.0835
$ws.Range("L76").Value = 3335.7144
$ws.Range("M76").Value = -6835.0835
$ws.Range("N76").Value = -3965.7144

$ws.Range("H79").Value = 6288.7744
$ws.Range("I79").Value = 7150.0835
$ws.Range("J79").Value = 3335.7144
$ws.Range("K79").Value = 7150.0835
$ws.Range("L79").Value = 3335.7144
$ws.Range("M79").Value = -6058.0835
$ws.Range("N79").Value = -5519.7144

$ws.Range("H137").Value = 41669144
$ws.Range("I137").Value = 1588.7333
$ws.Range("J137").Value = 111115064
$ws.Range("K137").Value = 4766.199900000001
$ws.Range("L137").Value = 333345192
$ws.Range("M137").Value = -2216.199900000001
$ws.Range("N137").Value = -333350292

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2561.6
$ws.Range("I32").Value = 2547.0708
$ws.Range("J32").Value = 4000
$ws.Range("K32").Value = 2547.0708
$ws.Range("L32").Value = 4000
$ws.Range("M32").Value = -2260.0708
$ws.Range("N32").Value = -4574

$ws.Range("H45").Value = 1018.84375
$ws.Range("I45").Value = 773.0476
$ws.Range("J45").Value = 1488.091
$ws.Range("K45").Value = 773.0476
$ws.Range("L45").Value = 1488.091
$ws.Range("M45").Value = -396.0476
$ws.Range("N45").Value = -2242.091

$ws.Range("H132").Value = 23760.848
$ws.Range("I132").Value = 2123.3901
$ws.Range("J132").Value = 201188
$ws.Range("K132").Value = 6370.1703
$ws.Range("L132").Value = 603564
$ws.Range("M132").Value = -3840.1703
$ws.Range("N132").Value = -608624

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1086.7778
$ws.Range("I99").Value = 978.3333
$ws.Range("J99").Value = 1303.6666
$ws.Range("K99").Value = 978.3333
$ws.Range("L99").Value = 1303.6666
$ws.Range("M99").Value = 519.6667
$ws.Range("N99").Value = -4299.6666

$ws.Range("H134").Value = 51655.227
$ws.Range("I134").Value = 70345.81
$ws.Range("J134").Value = 1813.6666
$ws.Range("K134").Value = 211037.43
$ws.Range("L134").Value = 5440.9998
$ws.Range("M134").Value = -208502.43
$ws.Range("N134").Value = -10510.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1298.3658
$ws.Range("I31").Value = 921.3939
$ws.Range("J31").Value = 2853.375
$ws.Range("K31").Value = 921.3939
$ws.Range("L31").Value = 2853.375
$ws.Range("M31").Value = -626.3939
$ws.Range("N31").Value = -3443.375

$ws.Range("H34").Value = 1298.3658
$ws.Range("I34").Value = 921.3939
$ws.Range("J34").Value = 2853.375
$ws.Range("K34").Value = 921.3939
$ws.Range("L34").Value = 2853.375
$ws.Range("M34").Value = -719.3939
$ws.Range("N34").Value = -3257.375

$ws.Range("H86").Value = 4283
$ws.Range("I86").Value = 4058.6667
$ws.Range("J86").Value = 4571.4287
$ws.Range("K86").Value = 4058.6667
$ws.Range("L86").Value = 4571.4287
$ws.Range("M86").Value = -2935.6667
$ws.Range("N86").Value = -6817.4287

$ws.Range("H89").Value = 4283
$ws.Range("I89").Value = 4058.6667
$ws.Range("J89").Value = 4571.4287
$ws.Range("K89").Value = 20293.3335
$ws.Range("L89").Value = 22857.1435
$ws.Range("M89").Value = -14677.3335
$ws.Range("N89").Value = -34089.14350000001

$ws.Range("H94").Value = 1698.8334
$ws.Range("I94").Value = 1199.6666
$ws.Range("J94").Value = 2198
$ws.Range("K94").Value = 1199.6666
$ws.Range("L94").Value = 2198
$ws.Range("M94").Value = -748.6666
$ws.Range("N94").Value = -3100

$ws.Range("H122").Value = 6387.3687
$ws.Range("I122").Value = 6653.3335
$ws.Range("J122").Value = 1600
$ws.Range("K122").Value = 19960.0005
$ws.Range("L122").Value = 4800
$ws.Range("M122").Value = -17510.0005
$ws.Range("N122").Value = -9700

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 953.26
$ws.Range("I5").Value = 839.3333
$ws.Range("J5").Value = 978.2683
$ws.Range("K5").Value = 2517.9999
$ws.Range("L5").Value = 2934.8049
$ws.Range("M5").Value = -2405.9999
$ws.Range("N5").Value = -3158.8049

$ws.Range("H131").Value = 2339.7273
$ws.Range("I131").Value = 17110
$ws.Range("J131").Value = 1636.381
$ws.Range("K131").Value = 51330
$ws.Range("L131").Value = 4909.143
$ws.Range("M131").Value = -46290
$ws.Range("N131").Value = -14989.143

$ws.Range("H135").Value = 953.26
$ws.Range("I135").Value = 839.3333
$ws.Range("J135").Value = 978.2683
$ws.Range("K135").Value = 7553.9997
$ws.Range("L135").Value = 8804.414699999999
$ws.Range("M135").Value = -5018.9997
$ws.Range("N135").Value = -13874.4147

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7358.8823
$ws.Range("I70").Value = 9653.5
$ws.Range("J70").Value = 4777.4375
$ws.Range("K70").Value = 9653.5
$ws.Range("L70").Value = 4777.4375
$ws.Range("M70").Value = -9383.5
$ws.Range("N70").Value = -5317.4375

$ws.Range("H73").Value = 7358.8823
$ws.Range("I73").Value = 9653.5
$ws.Range("J73").Value = 4777.4375
$ws.Range("K73").Value = 9653.5
$ws.Range("L73").Value = 4777.4375
$ws.Range("M73").Value = -8717.5
$ws.Range("N73").Value = -6649.4375

$ws.Range("H97").Value = 2331.4443
$ws.Range("I97").Value = 1991.6666
$ws.Range("J97").Value = 3011
$ws.Range("K97").Value = 1991.6666
$ws.Range("L97").Value = 3011
$ws.Range("M97").Value = -1495.6666
$ws.Range("N97").Value = -4003

$ws.Range("H107").Value = 399.56668
$ws.Range("I107").Value = 341.90475
$ws.Range("J107").Value = 534.1111
$ws.Range("K107").Value = 341.90475
$ws.Range("L107").Value = 534.1111
$ws.Range("M107").Value = 1578.09525
$ws.Range("N107").Value = -4374.1111

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1193.3334
$ws.Range("I7").Value = 1207.25
$ws.Range("J7").Value = 1123.75
$ws.Range("K7").Value = 1207.25
$ws.Range("L7").Value = 1123.75
$ws.Range("M7").Value = -1095.25
$ws.Range("N7").Value = -1347.75

$ws.Range("H126").Value = 1193.3334
$ws.Range("I126").Value = 1207.25
$ws.Range("J126").Value = 1123.75
$ws.Range("K126").Value = 3621.75
$ws.Range("L126").Value = 3371.25
$ws.Range("M126").Value = -1151.75
$ws.Range("N126").Value = -8311.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2359.4814
$ws.Range("I126").Value = 2027.579
$ws.Range("J126").Value = 3147.75
$ws.Range("K126").Value = 6082.737
$ws.Range("L126").Value = 9443.25
$ws.Range("M126").Value = -3612.737
$ws.Range("N126").Value = -14383.25

$ws.Range("H132").Value = 4093.3333
$ws.Range("I132").Value = 5519.0884
$ws.Range("J132").Value = 1241.8235
$ws.Range("K132").Value = 16557.2652
$ws.Range("L132").Value = 3725.4705
$ws.Range("M132").Value = -14027.2652
$ws.Range("N132").Value = -8785.470499999999
